# "Updated my status and the Project Plan"
#
# Weekly status workbook has two sheets: "7-2-13" (last week, unchanged
# content) and "7-9-13" (this week, gets refreshed with new status + the
# "Project Plan" meeting entry removed). We also flip the active tab over
# to the current week's sheet.

$wb = $excel.ActiveWorkbook
$wsLastWeek = $wb.Worksheets.Item(1)   # "7-2-13"
$wsThisWeek = $wb.Worksheets.Item(2)   # "7-9-13"

# --- This week's sheet ("7-9-13") gets the actual status update ---

# The "Total time:" literal count (F1) is stale until the week is filled
# in again - drop it.
$wsThisWeek.Range("F1").Clear()

# Row 4: first reading-log entry becomes the Android Tutorial, in progress.
$wsThisWeek.Range("A4").Value = "Android Tutorial"
$wsThisWeek.Range("B4").Value = 39997
$wsThisWeek.Range("C4").ClearContents()
$wsThisWeek.Range("D4").Value = 0.05
$wsThisWeek.Range("E4").Value = 0.25
$wsThisWeek.Range("F4").Clear()

# Row 5: the old Android Tutorial stub row is cleared back down to an
# empty placeholder (dates column keeps its date formatting, nothing else).
$wsThisWeek.Range("A5").Clear()
$wsThisWeek.Range("B5").ClearContents()
$wsThisWeek.Range("C5").ClearContents()
$wsThisWeek.Range("D5").Clear()
$wsThisWeek.Range("E5").Clear()

# Row 6: the "Project Plan Assistance" meeting entry is removed entirely.
$wsThisWeek.Range("A6:E6").Clear()

# Make "7-9-13" the active sheet/tab, with F2 selected.
$wsThisWeek.Activate()
$wsThisWeek.Range("F2").Select()
